# ============================================================================
# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board snapshot (currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and the derived LeveProfitNQ/HQ
# columns, H/I/J/K/L/M/N) for a handful of leves across the ALC, ARM, BSM,
# CRP, CUL, GSM and WVR sheets. A few rows had an NQ or HQ profit cell that
# is no longer meaningful (price now equals the leve cost) -- those cells are
# cleared outright rather than left with a stale value, and two CUL rows that
# previously had no HQ-profit (M) cell gain one.
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 108
$ws.Range("H108").Value = 79961.28999999999
$ws.Range("J108").Value = 79961.28999999999
$ws.Range("L108").Value = 79961.28999999999
$ws.Range("N108").Value = -87641.28999999999

# Row 109
$ws.Range("H109").Value = 49272.11
$ws.Range("J109").Value = 49272.11
$ws.Range("L109").Value = 49272.11
$ws.Range("N109").Value = -52046.11

# Row 110
$ws.Range("H110").Value = 67903.336
$ws.Range("J110").Value = 67903.336
$ws.Range("L110").Value = 67903.336
$ws.Range("N110").Value = -76083.336

# Row 120
$ws.Range("H120").Value = 49592.8
$ws.Range("J120").Value = 49592.8
$ws.Range("L120").Value = 49592.8
$ws.Range("N120").Value = -59268.8

# Row 123
$ws.Range("H123").Value = 68912.414
$ws.Range("J123").Value = 68912.414
$ws.Range("L123").Value = 68912.414
$ws.Range("N123").Value = -78712.414

# Row 129
$ws.Range("H129").Value = 2054.4443
$ws.Range("I129").Value = 2012.5
$ws.Range("K129").Value = 6037.5
$ws.Range("M129").Value = -1037.5

# Row 133
$ws.Range("H133").Value = 91459.375
$ws.Range("J133").Value = 91459.375
$ws.Range("L133").Value = 91459.375
$ws.Range("N133").Value = -101579.375

# Row 136
$ws.Range("H136").Value = 70162.164
$ws.Range("J136").Value = 70162.164
$ws.Range("L136").Value = 70162.164
$ws.Range("N136").Value = -80362.164

# Row 138
$ws.Range("H138").Value = 2116.4546
$ws.Range("I138").Value = 1828.1842
$ws.Range("J138").Value = 2760.8235
$ws.Range("K138").Value = 5484.5526
$ws.Range("L138").Value = 8282.470499999999
$ws.Range("M138").Value = -344.5526
$ws.Range("N138").Value = -18562.4705

# ---------------------------------------------------------------------------
# Sheet: ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 104
$ws.Range("H104").Value = 39897.6
$ws.Range("J104").Value = 39897.6
$ws.Range("L104").Value = 39897.6
$ws.Range("N104").Value = -46885.6

# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# Row 118
$ws.Range("H118").Value = 46426.285
$ws.Range("J118").Value = 46426.285
$ws.Range("L118").Value = 46426.285
$ws.Range("N118").Value = -49740.285

# ---------------------------------------------------------------------------
# Sheet: BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 13
$ws.Range("H13").Value = 86988
$ws.Range("J13").Value = 86988
$ws.Range("L13").Value = 86988
$ws.Range("N13").Value = -87324

# Row 53
$ws.Range("H53").Value = 68893
$ws.Range("J53").Value = 68893
$ws.Range("L53").Value = 68893
$ws.Range("N53").Value = -70041

# Row 108
$ws.Range("H108").Value = 99990
$ws.Range("J108").Value = 99990
$ws.Range("L108").Value = 99990
$ws.Range("N108").Value = -107670

# Row 109
$ws.Range("H109").Value = 99989.2
$ws.Range("J109").Value = 99989.2
$ws.Range("L109").Value = 99989.2
$ws.Range("N109").Value = -102763.2

# Row 110
$ws.Range("H110").Value = 69851.28999999999
$ws.Range("J110").Value = 69851.28999999999
$ws.Range("L110").Value = 69851.28999999999
$ws.Range("N110").Value = -78031.28999999999

# Row 114
$ws.Range("H114").Value = 89989.2
$ws.Range("J114").Value = 89989.2
$ws.Range("L114").Value = 89989.2
$ws.Range("N114").Value = -98667.2

# Row 118
$ws.Range("H118").Value = 72033.71000000001
$ws.Range("J118").Value = 74781.664
$ws.Range("L118").Value = 74781.664
$ws.Range("N118").Value = -78095.664

# Row 122
$ws.Range("H122").Value = 77773.164
$ws.Range("J122").Value = 77773.164
$ws.Range("L122").Value = 77773.164
$ws.Range("N122").Value = -87573.164

# Row 132
$ws.Range("H132").Value = 93282
$ws.Range("J132").Value = 93282
$ws.Range("L132").Value = 93282
$ws.Range("N132").Value = -103402

# Row 135
$ws.Range("H135").Value = 55785.875
$ws.Range("J135").Value = 55785.875
$ws.Range("L135").Value = 55785.875
$ws.Range("N135").Value = -65925.875

# Row 138
$ws.Range("H138").Value = 83886.75
$ws.Range("J138").Value = 83886.75
$ws.Range("L138").Value = 83886.75
$ws.Range("N138").Value = -94166.75

# Row 140
$ws.Range("H140").Value = 43576.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 43576.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 43576.5
$ws.Range("N140").Value = -53936.5
$ws.Range("M140").ClearContents()

# ---------------------------------------------------------------------------
# Sheet: CRP
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 9
$ws.Range("H9").Value = 30820.715
$ws.Range("J9").Value = 30820.715
$ws.Range("L9").Value = 30820.715
$ws.Range("N9").Value = -31156.715

# Row 31
$ws.Range("H31").Value = 2494.5881
$ws.Range("I31").Value = 2012.4166
$ws.Range("J31").Value = 3651.8
$ws.Range("K31").Value = 2012.4166
$ws.Range("L31").Value = 3651.8
$ws.Range("M31").Value = -1717.4166
$ws.Range("N31").Value = -4241.8

# Row 34
$ws.Range("H34").Value = 2494.5881
$ws.Range("I34").Value = 2012.4166
$ws.Range("J34").Value = 3651.8
$ws.Range("K34").Value = 2012.4166
$ws.Range("L34").Value = 3651.8
$ws.Range("M34").Value = -1810.4166
$ws.Range("N34").Value = -4055.8

# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()

# Row 114
$ws.Range("H114").Value = 45525.57
$ws.Range("J114").Value = 45525.57
$ws.Range("L114").Value = 45525.57
$ws.Range("N114").Value = -54203.57

# Row 116
$ws.Range("H116").Value = 45156
$ws.Range("J116").Value = 45156
$ws.Range("L116").Value = 45156
$ws.Range("N116").Value = -54334

# Row 119
$ws.Range("H119").Value = 63879.734
$ws.Range("J119").Value = 63879.734
$ws.Range("L119").Value = 63879.734
$ws.Range("N119").Value = -73555.734

# Row 138
$ws.Range("H138").Value = 53083.2
$ws.Range("J138").Value = 49926.75
$ws.Range("L138").Value = 49926.75
$ws.Range("N138").Value = -60206.75

# ---------------------------------------------------------------------------
# Sheet: CUL
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 907.0769
$ws.Range("I5").Value = 777
$ws.Range("J5").Value = 1199.75
$ws.Range("K5").Value = 2331
$ws.Range("L5").Value = 3599.25
$ws.Range("M5").Value = -2219
$ws.Range("N5").Value = -3823.25

# Row 63
$ws.Range("H63").Value = 200
$ws.Range("I63").Value = 200
$ws.Range("K63").Value = 600
$ws.Range("M63").Value = 149

# Row 66
$ws.Range("H66").Value = 200
$ws.Range("I66").Value = 200
$ws.Range("K66").Value = 1800
$ws.Range("M66").Value = 1944

# Row 68
$ws.Range("H68").Value = 10280.625
$ws.Range("I68").Value = 4374
$ws.Range("J68").Value = 12249.5
$ws.Range("K68").Value = 13122
$ws.Range("L68").Value = 36748.5
$ws.Range("M68").Value = -12311
$ws.Range("N68").Value = -38370.5

# Row 71
$ws.Range("H71").Value = 10280.625
$ws.Range("I71").Value = 4374
$ws.Range("J71").Value = 12249.5
$ws.Range("K71").Value = 39366
$ws.Range("L71").Value = 110245.5
$ws.Range("M71").Value = -35310
$ws.Range("N71").Value = -118357.5

# Row 107
$ws.Range("H107").Value = 369
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# Row 113
$ws.Range("H113").Value = 68038.47
$ws.Range("I113").Value = 1271.2858
$ws.Range("K113").Value = 3813.8574
$ws.Range("M113").Value = -1643.8574

# Row 135
$ws.Range("H135").Value = 907.0769
$ws.Range("I135").Value = 777
$ws.Range("J135").Value = 1199.75
$ws.Range("K135").Value = 6993
$ws.Range("L135").Value = 10797.75
$ws.Range("M135").Value = -4458
$ws.Range("N135").Value = -15867.75

# ---------------------------------------------------------------------------
# Sheet: GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 114
$ws.Range("H114").Value = 74087.336
$ws.Range("J114").Value = 74087.336
$ws.Range("L114").Value = 74087.336
$ws.Range("N114").Value = -82765.336

# Row 135
$ws.Range("H135").Value = 61385.777
$ws.Range("J135").Value = 61385.777
$ws.Range("L135").Value = 61385.777
$ws.Range("N135").Value = -71525.777

# Row 140
$ws.Range("H140").Value = 98496
$ws.Range("J140").Value = 98496
$ws.Range("L140").Value = 98496
$ws.Range("N140").Value = -108856

# ---------------------------------------------------------------------------
# Sheet: WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 28
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()

# Row 133
$ws.Range("H133").Value = 66637.60000000001
$ws.Range("J133").Value = 64422
$ws.Range("L133").Value = 64422
$ws.Range("N133").Value = -74542

# Row 136
$ws.Range("H136").Value = 1290.7778
$ws.Range("I136").Value = 737.2857
$ws.Range("K136").Value = 2211.8571
$ws.Range("M136").Value = 338.1428999999998
